$d = $word.ActiveDocument

function Set-ParagraphRuns($paragraph, $style, $words) {
    $runsXml = ""
    foreach ($w in $words) {
        $runsXml += '<w:r><w:t xml:space="preserve">' + $w + '</w:t></w:r>'
    }
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p><w:pPr><w:pStyle w:val="' + $style + '"/></w:pPr>' + $runsXml + '</w:p></w:body></w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
    $paragraph.Range.InsertXML($xml)
}

# Paragraph 1: Title
$p1 = $d.Paragraphs.Item(1)
Set-ParagraphRuns $p1 "Title" @("Answers:", " ", "Introduction", " ", "to", " ", "quadratic", " ", "equations")

# Paragraph 2: Author
$p2 = $d.Paragraphs.Item(2)
Set-ParagraphRuns $p2 "Author" @("Tom", " ", "Coleman")

# Paragraph 4: Abstract
$p4 = $d.Paragraphs.Item(4)
Set-ParagraphRuns $p4 "Abstract" @("Answers", " ", "to", " ", "questions", " ", "relating", " ", "to", " ", "the", " ", "guide", " ", "on", " ", "introduction", " ", "to", " ", "quadratic", " ", "equations.")
